$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting (style) for the brand-new rows 192 and 193 ---
# Column A uses the bold/bordered header-ish style, column E uses the datetime style.
$ws.Range("A191").Copy()
$ws.Range("A192:A193").PasteSpecial(-4122)
$ws.Range("E191").Copy()
$ws.Range("E192:E193").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 95 ---
$ws.Range("A95").Value = 93
$ws.Range("B95").Value = 6236612
$ws.Range("C95").Value = "Venezuela Primera Division"
$ws.Range("D95").Value = "Venezuela Primera Division"
$ws.Range("E95").Value = 45199.6875
$ws.Range("F95").Value = "Zamora"
$ws.Range("G95").Value = "Carabobo"
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 2
$ws.Range("J95").Value = "A"
$ws.Range("K95").Value = 3.2
$ws.Range("L95").Value = 3.1
$ws.Range("M95").Value = 2.15
$ws.Range("N95").Value = 4.5
$ws.Range("O95").Value = 3.3
$ws.Range("P95").Value = 1.75
$ws.Range("Q95").Value = 0.5
$ws.Range("R95").Value = 2
$ws.Range("S95").Value = 1.8
$ws.Range("T95").Value = 2.25
$ws.Range("U95").Value = 1.925
$ws.Range("V95").Value = 1.875
$ws.Range("W95").Value = -1
$ws.Range("X95").Value = -1
$ws.Range("Y95").Value = 0.75
$ws.Range("Z95").Value = -1
$ws.Range("AA95").Value = 0.8
$ws.Range("AB95").Value = -0.5
$ws.Range("AC95").Value = 0.4375

# --- Row 96 ---
$ws.Range("A96").Value = 94
$ws.Range("B96").Value = 6236252
$ws.Range("C96").Value = "Venezuela Primera Division"
$ws.Range("D96").Value = "Venezuela Primera Division"
$ws.Range("E96").Value = 45199.6875
$ws.Range("F96").Value = "Deportivo Tachira"
$ws.Range("G96").Value = "CD Hermanos Colmenares"
$ws.Range("H96").Value = 1
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = "H"
$ws.Range("K96").Value = 1.363
$ws.Range("L96").Value = 4.2
$ws.Range("M96").Value = 7.5
$ws.Range("N96").Value = 1.333
$ws.Range("O96").Value = 4.5
$ws.Range("P96").Value = 8
$ws.Range("Q96").Value = -1.5
$ws.Range("R96").Value = 2
$ws.Range("S96").Value = 1.8
$ws.Range("T96").Value = 2.5
$ws.Range("U96").Value = 1.925
$ws.Range("V96").Value = 1.875
$ws.Range("W96").Value = 0.333
$ws.Range("X96").Value = -1
$ws.Range("Y96").Value = -1
$ws.Range("Z96").Value = -1
$ws.Range("AA96").Value = 0.8
$ws.Range("AB96").Value = -1
$ws.Range("AC96").Value = 0.875

# --- Row 100 ---
$ws.Range("A100").Value = 98
$ws.Range("B100").Value = 6236614
$ws.Range("C100").Value = "Venezuela Primera Division"
$ws.Range("D100").Value = "Venezuela Primera Division"
$ws.Range("E100").Value = 45205.70833333334
$ws.Range("F100").Value = "Mineros"
$ws.Range("G100").Value = "Angostura FC"
$ws.Range("H100").Value = 1
$ws.Range("I100").Value = 2
$ws.Range("J100").Value = "A"
$ws.Range("K100").Value = 2.45
$ws.Range("L100").Value = 3.3
$ws.Range("M100").Value = 2.55
$ws.Range("N100").Value = 1.8
$ws.Range("O100").Value = 3.75
$ws.Range("P100").Value = 3.6
$ws.Range("Q100").Value = -0.5
$ws.Range("R100").Value = 1.825
$ws.Range("S100").Value = 1.975
$ws.Range("T100").Value = 2.75
$ws.Range("U100").Value = 1.8
$ws.Range("V100").Value = 2
$ws.Range("W100").Value = -1
$ws.Range("X100").Value = -1
$ws.Range("Y100").Value = 2.6
$ws.Range("Z100").Value = -1
$ws.Range("AA100").Value = 0.9750000000000001
$ws.Range("AB100").Value = 0.4
$ws.Range("AC100").Value = -0.5

# --- Row 101 ---
$ws.Range("A101").Value = 99
$ws.Range("B101").Value = 6236257
$ws.Range("C101").Value = "Venezuela Primera Division"
$ws.Range("D101").Value = "Venezuela Primera Division"
$ws.Range("E101").Value = 45205.70833333334
$ws.Range("F101").Value = "CD Hermanos Colmenares"
$ws.Range("G101").Value = "Zamora"
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 2
$ws.Range("J101").Value = "A"
$ws.Range("K101").Value = 2.3
$ws.Range("L101").Value = 3.2
$ws.Range("M101").Value = 2.8
$ws.Range("N101").Value = 1.666
$ws.Range("O101").Value = 3.8
$ws.Range("P101").Value = 4.2
$ws.Range("Q101").Value = -0.75
$ws.Range("R101").Value = 1.9
$ws.Range("S101").Value = 1.9
$ws.Range("T101").Value = 2.75
$ws.Range("U101").Value = 1.9
$ws.Range("V101").Value = 1.9
$ws.Range("W101").Value = -1
$ws.Range("X101").Value = -1
$ws.Range("Y101").Value = 3.2
$ws.Range("Z101").Value = -1
$ws.Range("AA101").Value = 0.8999999999999999
$ws.Range("AB101").Value = -1
$ws.Range("AC101").Value = 0.8999999999999999

# --- Row 162 ---
$ws.Range("A162").Value = 160
$ws.Range("B162").Value = 7952893
$ws.Range("C162").Value = "Venezuela Primera Division"
$ws.Range("D162").Value = "Venezuela Primera Division"
$ws.Range("E162").Value = 45366.83333333334
$ws.Range("F162").Value = "UCV"
$ws.Range("G162").Value = "Deportivo La Guaira"
$ws.Range("H162").Value = 1
$ws.Range("I162").Value = 1
$ws.Range("J162").Value = "D"
$ws.Range("K162").Value = 2.1
$ws.Range("L162").Value = 3
$ws.Range("M162").Value = 3.25
$ws.Range("N162").Value = 2.25
$ws.Range("O162").Value = 3.1
$ws.Range("P162").Value = 2.9
$ws.Range("Q162").Value = -0.25
$ws.Range("R162").Value = 2.025
$ws.Range("S162").Value = 1.775
$ws.Range("T162").Value = 2
$ws.Range("U162").Value = 1.8
$ws.Range("V162").Value = 2
$ws.Range("W162").Value = -1
$ws.Range("X162").Value = 2.1
$ws.Range("Y162").Value = -1
$ws.Range("Z162").Value = -0.5
$ws.Range("AA162").Value = 0.3875
$ws.Range("AB162").Value = 0
$ws.Range("AC162").Value = 0

# --- Row 163 ---
$ws.Range("A163").Value = 161
$ws.Range("B163").Value = 7952905
$ws.Range("C163").Value = "Venezuela Primera Division"
$ws.Range("D163").Value = "Venezuela Primera Division"
$ws.Range("E163").Value = 45366.83333333334
$ws.Range("F163").Value = "Angostura FC"
$ws.Range("G163").Value = "Deportivo Tachira"
$ws.Range("H163").Value = 2
$ws.Range("I163").Value = 0
$ws.Range("J163").Value = "H"
$ws.Range("K163").Value = 3.6
$ws.Range("L163").Value = 3.6
$ws.Range("M163").Value = 1.8
$ws.Range("N163").Value = 3.75
$ws.Range("O163").Value = 2.875
$ws.Range("P163").Value = 2.1
$ws.Range("Q163").Value = 0.25
$ws.Range("R163").Value = 1.95
$ws.Range("S163").Value = 1.85
$ws.Range("T163").Value = 2
$ws.Range("U163").Value = 2.025
$ws.Range("V163").Value = 1.775
$ws.Range("W163").Value = 2.75
$ws.Range("X163").Value = -1
$ws.Range("Y163").Value = -1
$ws.Range("Z163").Value = 0.95
$ws.Range("AA163").Value = -1
$ws.Range("AB163").Value = 0
$ws.Range("AC163").Value = 0

# --- Row 173 ---
$ws.Range("A173").Value = 171
$ws.Range("B173").Value = 7958192
$ws.Range("C173").Value = "Venezuela Primera Division"
$ws.Range("D173").Value = "Venezuela Primera Division"
$ws.Range("E173").Value = 45371.89583333334
$ws.Range("F173").Value = "Deportivo Tachira"
$ws.Range("G173").Value = "Monagas"
$ws.Range("H173").Value = 1
$ws.Range("I173").Value = 0
$ws.Range("J173").Value = "H"
$ws.Range("K173").Value = 1.666
$ws.Range("L173").Value = 3.4
$ws.Range("M173").Value = 4.5
$ws.Range("N173").Value = 1.95
$ws.Range("O173").Value = 3.25
$ws.Range("P173").Value = 3.5
$ws.Range("Q173").Value = -0.5
$ws.Range("R173").Value = 1.975
$ws.Range("S173").Value = 1.825
$ws.Range("T173").Value = 2.25
$ws.Range("U173").Value = 2.025
$ws.Range("V173").Value = 1.775
$ws.Range("W173").Value = 0.95
$ws.Range("X173").Value = -1
$ws.Range("Y173").Value = -1
$ws.Range("Z173").Value = 0.9750000000000001
$ws.Range("AA173").Value = -1
$ws.Range("AB173").Value = -1
$ws.Range("AC173").Value = 0.7749999999999999

# --- Row 174 ---
$ws.Range("A174").Value = 172
$ws.Range("B174").Value = 7958193
$ws.Range("C174").Value = "Venezuela Primera Division"
$ws.Range("D174").Value = "Venezuela Primera Division"
$ws.Range("E174").Value = 45371.89583333334
$ws.Range("F174").Value = "Zamora"
$ws.Range("G174").Value = "Academia Puerto Cabello"
$ws.Range("H174").Value = 0
$ws.Range("I174").Value = 0
$ws.Range("J174").Value = "D"
$ws.Range("K174").Value = 3.75
$ws.Range("L174").Value = 3.3
$ws.Range("M174").Value = 1.85
$ws.Range("N174").Value = 3.1
$ws.Range("O174").Value = 3.2
$ws.Range("P174").Value = 2.1
$ws.Range("Q174").Value = 0.25
$ws.Range("R174").Value = 1.875
$ws.Range("S174").Value = 1.925
$ws.Range("T174").Value = 2.25
$ws.Range("U174").Value = 2.025
$ws.Range("V174").Value = 1.775
$ws.Range("W174").Value = -1
$ws.Range("X174").Value = 2.2
$ws.Range("Y174").Value = -1
$ws.Range("Z174").Value = 0.4375
$ws.Range("AA174").Value = -0.5
$ws.Range("AB174").Value = -1
$ws.Range("AC174").Value = 0.7749999999999999

# --- Row 190 ---
$ws.Range("A190").Value = 188
$ws.Range("B190").Value = 8054935
$ws.Range("C190").Value = "Venezuela Primera Division"
$ws.Range("D190").Value = "Venezuela Primera Division"
$ws.Range("E190").Value = 45388.66666666666
$ws.Range("F190").Value = "Estudiantes Merida"
$ws.Range("G190").Value = "Zamora"
$ws.Range("K190").Value = 2.3
$ws.Range("L190").Value = 3.1
$ws.Range("M190").Value = 2.875
$ws.Range("N190").Value = 2.375
$ws.Range("O190").Value = 3.25
$ws.Range("P190").Value = 2.6
$ws.Range("Q190").Value = 0
$ws.Range("R190").Value = 1.8
$ws.Range("S190").Value = 2
$ws.Range("T190").Value = 2.5
$ws.Range("U190").Value = 1.925
$ws.Range("V190").Value = 1.875
$ws.Range("W190").Value = 0
$ws.Range("X190").Value = 0
$ws.Range("Y190").Value = 0
$ws.Range("Z190").Value = 0
$ws.Range("AA190").Value = 0

# --- Row 191 ---
$ws.Range("A191").Value = 189
$ws.Range("B191").Value = 8054936
$ws.Range("C191").Value = "Venezuela Primera Division"
$ws.Range("D191").Value = "Venezuela Primera Division"
$ws.Range("E191").Value = 45388.78125
$ws.Range("F191").Value = "Carabobo"
$ws.Range("G191").Value = "Deportivo Rayo Zuliano"
$ws.Range("K191").Value = 1.444
$ws.Range("L191").Value = 3.75
$ws.Range("M191").Value = 7
$ws.Range("N191").Value = 1.45
$ws.Range("O191").Value = 3.75
$ws.Range("P191").Value = 7
$ws.Range("Q191").Value = -1.25
$ws.Range("R191").Value = 2
$ws.Range("S191").Value = 1.8
$ws.Range("T191").Value = 2.5
$ws.Range("U191").Value = 2
$ws.Range("V191").Value = 1.8
$ws.Range("W191").Value = 0
$ws.Range("X191").Value = 0
$ws.Range("Y191").Value = 0
$ws.Range("Z191").Value = 0
$ws.Range("AA191").Value = 0

# --- Row 192 ---
$ws.Range("A192").Value = 190
$ws.Range("B192").Value = 8054937
$ws.Range("C192").Value = "Venezuela Primera Division"
$ws.Range("D192").Value = "Venezuela Primera Division"
$ws.Range("E192").Value = 45388.89583333334
$ws.Range("F192").Value = "Deportivo La Guaira"
$ws.Range("G192").Value = "Portuguesa"
$ws.Range("K192").Value = 2.25
$ws.Range("L192").Value = 2.875
$ws.Range("M192").Value = 3.25
$ws.Range("N192").Value = 2.1
$ws.Range("O192").Value = 2.9
$ws.Range("P192").Value = 3.6
$ws.Range("Q192").Value = -0.25
$ws.Range("R192").Value = 1.825
$ws.Range("S192").Value = 1.975
$ws.Range("T192").Value = 2
$ws.Range("U192").Value = 1.85
$ws.Range("V192").Value = 1.95
$ws.Range("W192").Value = 0
$ws.Range("X192").Value = 0
$ws.Range("Y192").Value = 0
$ws.Range("Z192").Value = 0
$ws.Range("AA192").Value = 0

# --- Row 193 ---
$ws.Range("A193").Value = 191
$ws.Range("B193").Value = 8054938
$ws.Range("C193").Value = "Venezuela Primera Division"
$ws.Range("D193").Value = "Venezuela Primera Division"
$ws.Range("E193").Value = 45389.70833333334
$ws.Range("F193").Value = "Metropolitanos FC"
$ws.Range("G193").Value = "Caracas"
$ws.Range("K193").Value = 3
$ws.Range("L193").Value = 2.875
$ws.Range("M193").Value = 2.375
$ws.Range("N193").Value = 2.5
$ws.Range("O193").Value = 2.875
$ws.Range("P193").Value = 2.75
$ws.Range("Q193").Value = 0
$ws.Range("R193").Value = 1.775
$ws.Range("S193").Value = 2.025
$ws.Range("T193").Value = 2.25
$ws.Range("U193").Value = 2.05
$ws.Range("V193").Value = 1.75
$ws.Range("W193").Value = 0
$ws.Range("X193").Value = 0
$ws.Range("Y193").Value = 0
$ws.Range("Z193").Value = 0
$ws.Range("AA193").Value = 0

